$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52, shifting existing rows 52-68 down to 53-69
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new weekly record
$ws.Range("A52").Value = 11
$ws.Range("B52").Value = "Vega Monumental Concepción"
$ws.Range("C52").Value = "Bíobío"
$ws.Range("D52").Value = 44588
$ws.Range("E52").Value = 8
$ws.Range("F52").Value = 100112021
$ws.Range("G52").Value = "Ají"
$ws.Range("H52").Value = "Chilena(o)"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 50
$ws.Range("K52").Value = 22000
$ws.Range("L52").Value = 23000
$ws.Range("M52").Value = 22400
$ws.Range("N52").Value = "$/caja 12 kilos"
$ws.Range("O52").Value = "Región Metropolitana"
$ws.Range("P52").Value = 1867
$ws.Range("Q52").Value = 12
$ws.Range("R52").Value = "Hortaliza"
